$wb = $excel.ActiveWorkbook

# --- PIR sheet: append new rows ---
$ws = $wb.Worksheets.Item("PIR")
# Row 163
$ws.Cells.Item(163, 1).NumberFormat = "@"
$ws.Cells.Item(163, 1).Value = '2026-01-28'
$ws.Cells.Item(163, 1).NumberFormat = "General"
$ws.Cells.Item(163, 2).Value = '18:14:59'
$ws.Cells.Item(163, 3).Value = '18:00'
$ws.Cells.Item(163, 4).Value = 'Bathroom'
$ws.Cells.Item(163, 5).Value = 'No Motion'
$ws.Cells.Item(163, 6).Value = 'Inactive'

# Row 164
$ws.Cells.Item(164, 1).NumberFormat = "@"
$ws.Cells.Item(164, 1).Value = '2026-01-28'
$ws.Cells.Item(164, 1).NumberFormat = "General"
$ws.Cells.Item(164, 2).Value = '18:15:04'
$ws.Cells.Item(164, 3).Value = '18:00'
$ws.Cells.Item(164, 4).Value = 'Bathroom'
$ws.Cells.Item(164, 5).Value = 'No Motion'
$ws.Cells.Item(164, 6).Value = 'Inactive'

# Row 165
$ws.Cells.Item(165, 1).NumberFormat = "@"
$ws.Cells.Item(165, 1).Value = '2026-01-28'
$ws.Cells.Item(165, 1).NumberFormat = "General"
$ws.Cells.Item(165, 2).Value = '18:15:06'
$ws.Cells.Item(165, 3).Value = '18:00'
$ws.Cells.Item(165, 4).Value = 'Bathroom'
$ws.Cells.Item(165, 5).Value = 'No Motion'
$ws.Cells.Item(165, 6).Value = 'Inactive'

# Row 166
$ws.Cells.Item(166, 1).NumberFormat = "@"
$ws.Cells.Item(166, 1).Value = '2026-01-28'
$ws.Cells.Item(166, 1).NumberFormat = "General"
$ws.Cells.Item(166, 2).Value = '18:15:07'
$ws.Cells.Item(166, 3).Value = '18:00'
$ws.Cells.Item(166, 4).Value = 'Bathroom'
$ws.Cells.Item(166, 5).Value = 'No Motion'
$ws.Cells.Item(166, 6).Value = 'Inactive'

# Row 167
$ws.Cells.Item(167, 1).NumberFormat = "@"
$ws.Cells.Item(167, 1).Value = '2026-01-28'
$ws.Cells.Item(167, 1).NumberFormat = "General"
$ws.Cells.Item(167, 2).Value = '18:15:11'
$ws.Cells.Item(167, 3).Value = '18:00'
$ws.Cells.Item(167, 4).Value = 'Bathroom'
$ws.Cells.Item(167, 5).Value = 'No Motion'
$ws.Cells.Item(167, 6).Value = 'Inactive'

# Row 168
$ws.Cells.Item(168, 1).NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = '2026-01-28'
$ws.Cells.Item(168, 1).NumberFormat = "General"
$ws.Cells.Item(168, 2).Value = '18:15:18'
$ws.Cells.Item(168, 3).Value = '18:00'
$ws.Cells.Item(168, 4).Value = 'Bathroom'
$ws.Cells.Item(168, 5).Value = 'No Motion'
$ws.Cells.Item(168, 6).Value = 'Inactive'

# Row 169
$ws.Cells.Item(169, 1).NumberFormat = "@"
$ws.Cells.Item(169, 1).Value = '2026-01-28'
$ws.Cells.Item(169, 1).NumberFormat = "General"
$ws.Cells.Item(169, 2).Value = '18:15:22'
$ws.Cells.Item(169, 3).Value = '18:00'
$ws.Cells.Item(169, 4).Value = 'Bathroom'
$ws.Cells.Item(169, 5).Value = 'No Motion'
$ws.Cells.Item(169, 6).Value = 'Inactive'

# Row 170
$ws.Cells.Item(170, 1).NumberFormat = "@"
$ws.Cells.Item(170, 1).Value = '2026-01-28'
$ws.Cells.Item(170, 1).NumberFormat = "General"
$ws.Cells.Item(170, 2).Value = '18:15:26'
$ws.Cells.Item(170, 3).Value = '18:00'
$ws.Cells.Item(170, 4).Value = 'Bathroom'
$ws.Cells.Item(170, 5).Value = 'No Motion'
$ws.Cells.Item(170, 6).Value = 'Inactive'

# Row 171
$ws.Cells.Item(171, 1).NumberFormat = "@"
$ws.Cells.Item(171, 1).Value = '2026-01-28'
$ws.Cells.Item(171, 1).NumberFormat = "General"
$ws.Cells.Item(171, 2).Value = '18:15:31'
$ws.Cells.Item(171, 3).Value = '18:00'
$ws.Cells.Item(171, 4).Value = 'Bathroom'
$ws.Cells.Item(171, 5).Value = 'No Motion'
$ws.Cells.Item(171, 6).Value = 'Inactive'

# Row 172
$ws.Cells.Item(172, 1).NumberFormat = "@"
$ws.Cells.Item(172, 1).Value = '2026-01-28'
$ws.Cells.Item(172, 1).NumberFormat = "General"
$ws.Cells.Item(172, 2).Value = '18:15:38'
$ws.Cells.Item(172, 3).Value = '18:00'
$ws.Cells.Item(172, 4).Value = 'Bathroom'
$ws.Cells.Item(172, 5).Value = 'No Motion'
$ws.Cells.Item(172, 6).Value = 'Inactive'

# Row 173
$ws.Cells.Item(173, 1).NumberFormat = "@"
$ws.Cells.Item(173, 1).Value = '2026-01-28'
$ws.Cells.Item(173, 1).NumberFormat = "General"
$ws.Cells.Item(173, 2).Value = '18:15:42'
$ws.Cells.Item(173, 3).Value = '18:00'
$ws.Cells.Item(173, 4).Value = 'Bathroom'
$ws.Cells.Item(173, 5).Value = 'No Motion'
$ws.Cells.Item(173, 6).Value = 'Inactive'

# Row 174
$ws.Cells.Item(174, 1).NumberFormat = "@"
$ws.Cells.Item(174, 1).Value = '2026-01-28'
$ws.Cells.Item(174, 1).NumberFormat = "General"
$ws.Cells.Item(174, 2).Value = '18:15:46'
$ws.Cells.Item(174, 3).Value = '18:00'
$ws.Cells.Item(174, 4).Value = 'Bathroom'
$ws.Cells.Item(174, 5).Value = 'No Motion'
$ws.Cells.Item(174, 6).Value = 'Inactive'

# Row 175
$ws.Cells.Item(175, 1).NumberFormat = "@"
$ws.Cells.Item(175, 1).Value = '2026-01-28'
$ws.Cells.Item(175, 1).NumberFormat = "General"
$ws.Cells.Item(175, 2).Value = '18:15:51'
$ws.Cells.Item(175, 3).Value = '18:00'
$ws.Cells.Item(175, 4).Value = 'Bathroom'
$ws.Cells.Item(175, 5).Value = 'No Motion'
$ws.Cells.Item(175, 6).Value = 'Inactive'

# Row 176
$ws.Cells.Item(176, 1).NumberFormat = "@"
$ws.Cells.Item(176, 1).Value = '2026-01-28'
$ws.Cells.Item(176, 1).NumberFormat = "General"
$ws.Cells.Item(176, 2).Value = '18:15:58'
$ws.Cells.Item(176, 3).Value = '18:00'
$ws.Cells.Item(176, 4).Value = 'Bathroom'
$ws.Cells.Item(176, 5).Value = 'No Motion'
$ws.Cells.Item(176, 6).Value = 'Inactive'

# --- Humidity sheet: append new rows ---
$ws = $wb.Worksheets.Item("Humidity")
# Row 157
$ws.Cells.Item(157, 1).NumberFormat = "@"
$ws.Cells.Item(157, 1).Value = '2026-01-28'
$ws.Cells.Item(157, 1).NumberFormat = "General"
$ws.Cells.Item(157, 2).Value = '18:15:00'
$ws.Cells.Item(157, 3).Value = '18:00'
$ws.Cells.Item(157, 4).Value = 'Bathroom'
$ws.Cells.Item(157, 5).NumberFormat = "@"
$ws.Cells.Item(157, 5).Value = '88.2%'
$ws.Cells.Item(157, 5).NumberFormat = "General"
$ws.Cells.Item(157, 6).Value = 'Active'

# Row 158
$ws.Cells.Item(158, 1).NumberFormat = "@"
$ws.Cells.Item(158, 1).Value = '2026-01-28'
$ws.Cells.Item(158, 1).NumberFormat = "General"
$ws.Cells.Item(158, 2).Value = '18:15:02'
$ws.Cells.Item(158, 3).Value = '18:00'
$ws.Cells.Item(158, 4).Value = 'Bathroom'
$ws.Cells.Item(158, 5).NumberFormat = "@"
$ws.Cells.Item(158, 5).Value = '88.2%'
$ws.Cells.Item(158, 5).NumberFormat = "General"
$ws.Cells.Item(158, 6).Value = 'Active'

# Row 159
$ws.Cells.Item(159, 1).NumberFormat = "@"
$ws.Cells.Item(159, 1).Value = '2026-01-28'
$ws.Cells.Item(159, 1).NumberFormat = "General"
$ws.Cells.Item(159, 2).Value = '18:15:05'
$ws.Cells.Item(159, 3).Value = '18:00'
$ws.Cells.Item(159, 4).Value = 'Bathroom'
$ws.Cells.Item(159, 5).NumberFormat = "@"
$ws.Cells.Item(159, 5).Value = '87.3%'
$ws.Cells.Item(159, 5).NumberFormat = "General"
$ws.Cells.Item(159, 6).Value = 'Active'

# Row 160
$ws.Cells.Item(160, 1).NumberFormat = "@"
$ws.Cells.Item(160, 1).Value = '2026-01-28'
$ws.Cells.Item(160, 1).NumberFormat = "General"
$ws.Cells.Item(160, 2).Value = '18:15:08'
$ws.Cells.Item(160, 3).Value = '18:00'
$ws.Cells.Item(160, 4).Value = 'Bathroom'
$ws.Cells.Item(160, 5).NumberFormat = "@"
$ws.Cells.Item(160, 5).Value = '87.3%'
$ws.Cells.Item(160, 5).NumberFormat = "General"
$ws.Cells.Item(160, 6).Value = 'Active'

# Row 161
$ws.Cells.Item(161, 1).NumberFormat = "@"
$ws.Cells.Item(161, 1).Value = '2026-01-28'
$ws.Cells.Item(161, 1).NumberFormat = "General"
$ws.Cells.Item(161, 2).Value = '18:15:12'
$ws.Cells.Item(161, 3).Value = '18:00'
$ws.Cells.Item(161, 4).Value = 'Bathroom'
$ws.Cells.Item(161, 5).NumberFormat = "@"
$ws.Cells.Item(161, 5).Value = '88.2%'
$ws.Cells.Item(161, 5).NumberFormat = "General"
$ws.Cells.Item(161, 6).Value = 'Active'

# Row 162
$ws.Cells.Item(162, 1).NumberFormat = "@"
$ws.Cells.Item(162, 1).Value = '2026-01-28'
$ws.Cells.Item(162, 1).NumberFormat = "General"
$ws.Cells.Item(162, 2).Value = '18:15:16'
$ws.Cells.Item(162, 3).Value = '18:00'
$ws.Cells.Item(162, 4).Value = 'Bathroom'
$ws.Cells.Item(162, 5).NumberFormat = "@"
$ws.Cells.Item(162, 5).Value = '88.2%'
$ws.Cells.Item(162, 5).NumberFormat = "General"
$ws.Cells.Item(162, 6).Value = 'Active'

# Row 163
$ws.Cells.Item(163, 1).NumberFormat = "@"
$ws.Cells.Item(163, 1).Value = '2026-01-28'
$ws.Cells.Item(163, 1).NumberFormat = "General"
$ws.Cells.Item(163, 2).Value = '18:15:20'
$ws.Cells.Item(163, 3).Value = '18:00'
$ws.Cells.Item(163, 4).Value = 'Bathroom'
$ws.Cells.Item(163, 5).NumberFormat = "@"
$ws.Cells.Item(163, 5).Value = '87.3%'
$ws.Cells.Item(163, 5).NumberFormat = "General"
$ws.Cells.Item(163, 6).Value = 'Active'

# Row 164
$ws.Cells.Item(164, 1).NumberFormat = "@"
$ws.Cells.Item(164, 1).Value = '2026-01-28'
$ws.Cells.Item(164, 1).NumberFormat = "General"
$ws.Cells.Item(164, 2).Value = '18:15:24'
$ws.Cells.Item(164, 3).Value = '18:00'
$ws.Cells.Item(164, 4).Value = 'Bathroom'
$ws.Cells.Item(164, 5).NumberFormat = "@"
$ws.Cells.Item(164, 5).Value = '88.2%'
$ws.Cells.Item(164, 5).NumberFormat = "General"
$ws.Cells.Item(164, 6).Value = 'Active'

# Row 165
$ws.Cells.Item(165, 1).NumberFormat = "@"
$ws.Cells.Item(165, 1).Value = '2026-01-28'
$ws.Cells.Item(165, 1).NumberFormat = "General"
$ws.Cells.Item(165, 2).Value = '18:15:32'
$ws.Cells.Item(165, 3).Value = '18:00'
$ws.Cells.Item(165, 4).Value = 'Bathroom'
$ws.Cells.Item(165, 5).NumberFormat = "@"
$ws.Cells.Item(165, 5).Value = '87.3%'
$ws.Cells.Item(165, 5).NumberFormat = "General"
$ws.Cells.Item(165, 6).Value = 'Active'

# Row 166
$ws.Cells.Item(166, 1).NumberFormat = "@"
$ws.Cells.Item(166, 1).Value = '2026-01-28'
$ws.Cells.Item(166, 1).NumberFormat = "General"
$ws.Cells.Item(166, 2).Value = '18:15:36'
$ws.Cells.Item(166, 3).Value = '18:00'
$ws.Cells.Item(166, 4).Value = 'Bathroom'
$ws.Cells.Item(166, 5).NumberFormat = "@"
$ws.Cells.Item(166, 5).Value = '88.2%'
$ws.Cells.Item(166, 5).NumberFormat = "General"
$ws.Cells.Item(166, 6).Value = 'Active'

# Row 167
$ws.Cells.Item(167, 1).NumberFormat = "@"
$ws.Cells.Item(167, 1).Value = '2026-01-28'
$ws.Cells.Item(167, 1).NumberFormat = "General"
$ws.Cells.Item(167, 2).Value = '18:15:40'
$ws.Cells.Item(167, 3).Value = '18:00'
$ws.Cells.Item(167, 4).Value = 'Bathroom'
$ws.Cells.Item(167, 5).NumberFormat = "@"
$ws.Cells.Item(167, 5).Value = '87.3%'
$ws.Cells.Item(167, 5).NumberFormat = "General"
$ws.Cells.Item(167, 6).Value = 'Active'

# Row 168
$ws.Cells.Item(168, 1).NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = '2026-01-28'
$ws.Cells.Item(168, 1).NumberFormat = "General"
$ws.Cells.Item(168, 2).Value = '18:15:44'
$ws.Cells.Item(168, 3).Value = '18:00'
$ws.Cells.Item(168, 4).Value = 'Bathroom'
$ws.Cells.Item(168, 5).NumberFormat = "@"
$ws.Cells.Item(168, 5).Value = '88.2%'
$ws.Cells.Item(168, 5).NumberFormat = "General"
$ws.Cells.Item(168, 6).Value = 'Active'

# Row 169
$ws.Cells.Item(169, 1).NumberFormat = "@"
$ws.Cells.Item(169, 1).Value = '2026-01-28'
$ws.Cells.Item(169, 1).NumberFormat = "General"
$ws.Cells.Item(169, 2).Value = '18:15:52'
$ws.Cells.Item(169, 3).Value = '18:00'
$ws.Cells.Item(169, 4).Value = 'Bathroom'
$ws.Cells.Item(169, 5).NumberFormat = "@"
$ws.Cells.Item(169, 5).Value = '88.2%'
$ws.Cells.Item(169, 5).NumberFormat = "General"
$ws.Cells.Item(169, 6).Value = 'Active'

# Row 170
$ws.Cells.Item(170, 1).NumberFormat = "@"
$ws.Cells.Item(170, 1).Value = '2026-01-28'
$ws.Cells.Item(170, 1).NumberFormat = "General"
$ws.Cells.Item(170, 2).Value = '18:15:56'
$ws.Cells.Item(170, 3).Value = '18:00'
$ws.Cells.Item(170, 4).Value = 'Bathroom'
$ws.Cells.Item(170, 5).NumberFormat = "@"
$ws.Cells.Item(170, 5).Value = '88.3%'
$ws.Cells.Item(170, 5).NumberFormat = "General"
$ws.Cells.Item(170, 6).Value = 'Active'

# --- Temperature sheet: append new rows ---
$ws = $wb.Worksheets.Item("Temperature")
# Row 157
$ws.Cells.Item(157, 1).NumberFormat = "@"
$ws.Cells.Item(157, 1).Value = '2026-01-28'
$ws.Cells.Item(157, 1).NumberFormat = "General"
$ws.Cells.Item(157, 2).Value = '18:15:01'
$ws.Cells.Item(157, 3).Value = '18:00'
$ws.Cells.Item(157, 4).Value = 'Bathroom'
$ws.Cells.Item(157, 5).Value = '22.9C'
$ws.Cells.Item(157, 6).Value = 'Active'

# Row 158
$ws.Cells.Item(158, 1).NumberFormat = "@"
$ws.Cells.Item(158, 1).Value = '2026-01-28'
$ws.Cells.Item(158, 1).NumberFormat = "General"
$ws.Cells.Item(158, 2).Value = '18:15:03'
$ws.Cells.Item(158, 3).Value = '18:00'
$ws.Cells.Item(158, 4).Value = 'Bathroom'
$ws.Cells.Item(158, 5).Value = '22.9C'
$ws.Cells.Item(158, 6).Value = 'Active'

# Row 159
$ws.Cells.Item(159, 1).NumberFormat = "@"
$ws.Cells.Item(159, 1).Value = '2026-01-28'
$ws.Cells.Item(159, 1).NumberFormat = "General"
$ws.Cells.Item(159, 2).Value = '18:15:06'
$ws.Cells.Item(159, 3).Value = '18:00'
$ws.Cells.Item(159, 4).Value = 'Bathroom'
$ws.Cells.Item(159, 5).Value = '22.9C'
$ws.Cells.Item(159, 6).Value = 'Active'

# Row 160
$ws.Cells.Item(160, 1).NumberFormat = "@"
$ws.Cells.Item(160, 1).Value = '2026-01-28'
$ws.Cells.Item(160, 1).NumberFormat = "General"
$ws.Cells.Item(160, 2).Value = '18:15:09'
$ws.Cells.Item(160, 3).Value = '18:00'
$ws.Cells.Item(160, 4).Value = 'Bathroom'
$ws.Cells.Item(160, 5).Value = '22.9C'
$ws.Cells.Item(160, 6).Value = 'Active'

# Row 161
$ws.Cells.Item(161, 1).NumberFormat = "@"
$ws.Cells.Item(161, 1).Value = '2026-01-28'
$ws.Cells.Item(161, 1).NumberFormat = "General"
$ws.Cells.Item(161, 2).Value = '18:15:13'
$ws.Cells.Item(161, 3).Value = '18:00'
$ws.Cells.Item(161, 4).Value = 'Bathroom'
$ws.Cells.Item(161, 5).Value = '22.9C'
$ws.Cells.Item(161, 6).Value = 'Active'

# Row 162
$ws.Cells.Item(162, 1).NumberFormat = "@"
$ws.Cells.Item(162, 1).Value = '2026-01-28'
$ws.Cells.Item(162, 1).NumberFormat = "General"
$ws.Cells.Item(162, 2).Value = '18:15:17'
$ws.Cells.Item(162, 3).Value = '18:00'
$ws.Cells.Item(162, 4).Value = 'Bathroom'
$ws.Cells.Item(162, 5).Value = '22.9C'
$ws.Cells.Item(162, 6).Value = 'Active'

# Row 163
$ws.Cells.Item(163, 1).NumberFormat = "@"
$ws.Cells.Item(163, 1).Value = '2026-01-28'
$ws.Cells.Item(163, 1).NumberFormat = "General"
$ws.Cells.Item(163, 2).Value = '18:15:21'
$ws.Cells.Item(163, 3).Value = '18:00'
$ws.Cells.Item(163, 4).Value = 'Bathroom'
$ws.Cells.Item(163, 5).Value = '22.9C'
$ws.Cells.Item(163, 6).Value = 'Active'

# Row 164
$ws.Cells.Item(164, 1).NumberFormat = "@"
$ws.Cells.Item(164, 1).Value = '2026-01-28'
$ws.Cells.Item(164, 1).NumberFormat = "General"
$ws.Cells.Item(164, 2).Value = '18:15:25'
$ws.Cells.Item(164, 3).Value = '18:00'
$ws.Cells.Item(164, 4).Value = 'Bathroom'
$ws.Cells.Item(164, 5).Value = '22.9C'
$ws.Cells.Item(164, 6).Value = 'Active'

# Row 165
$ws.Cells.Item(165, 1).NumberFormat = "@"
$ws.Cells.Item(165, 1).Value = '2026-01-28'
$ws.Cells.Item(165, 1).NumberFormat = "General"
$ws.Cells.Item(165, 2).Value = '18:15:33'
$ws.Cells.Item(165, 3).Value = '18:00'
$ws.Cells.Item(165, 4).Value = 'Bathroom'
$ws.Cells.Item(165, 5).Value = '22.9C'
$ws.Cells.Item(165, 6).Value = 'Active'

# Row 166
$ws.Cells.Item(166, 1).NumberFormat = "@"
$ws.Cells.Item(166, 1).Value = '2026-01-28'
$ws.Cells.Item(166, 1).NumberFormat = "General"
$ws.Cells.Item(166, 2).Value = '18:15:37'
$ws.Cells.Item(166, 3).Value = '18:00'
$ws.Cells.Item(166, 4).Value = 'Bathroom'
$ws.Cells.Item(166, 5).Value = '22.9C'
$ws.Cells.Item(166, 6).Value = 'Active'

# Row 167
$ws.Cells.Item(167, 1).NumberFormat = "@"
$ws.Cells.Item(167, 1).Value = '2026-01-28'
$ws.Cells.Item(167, 1).NumberFormat = "General"
$ws.Cells.Item(167, 2).Value = '18:15:41'
$ws.Cells.Item(167, 3).Value = '18:00'
$ws.Cells.Item(167, 4).Value = 'Bathroom'
$ws.Cells.Item(167, 5).Value = '22.9C'
$ws.Cells.Item(167, 6).Value = 'Active'

# Row 168
$ws.Cells.Item(168, 1).NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = '2026-01-28'
$ws.Cells.Item(168, 1).NumberFormat = "General"
$ws.Cells.Item(168, 2).Value = '18:15:45'
$ws.Cells.Item(168, 3).Value = '18:00'
$ws.Cells.Item(168, 4).Value = 'Bathroom'
$ws.Cells.Item(168, 5).Value = '22.9C'
$ws.Cells.Item(168, 6).Value = 'Active'

# Row 169
$ws.Cells.Item(169, 1).NumberFormat = "@"
$ws.Cells.Item(169, 1).Value = '2026-01-28'
$ws.Cells.Item(169, 1).NumberFormat = "General"
$ws.Cells.Item(169, 2).Value = '18:15:53'
$ws.Cells.Item(169, 3).Value = '18:00'
$ws.Cells.Item(169, 4).Value = 'Bathroom'
$ws.Cells.Item(169, 5).Value = '22.9C'
$ws.Cells.Item(169, 6).Value = 'Active'

# Row 170
$ws.Cells.Item(170, 1).NumberFormat = "@"
$ws.Cells.Item(170, 1).Value = '2026-01-28'
$ws.Cells.Item(170, 1).NumberFormat = "General"
$ws.Cells.Item(170, 2).Value = '18:15:57'
$ws.Cells.Item(170, 3).Value = '18:00'
$ws.Cells.Item(170, 4).Value = 'Bathroom'
$ws.Cells.Item(170, 5).Value = '23.0C'
$ws.Cells.Item(170, 6).Value = 'Active'
